# 14 Nov (2nd commit)
# Adds a new "CreateContacts" worksheet (with a small firstname/lastname/
# companyname/runmode table), flips the existing test_suite runmode flags
# for the three pre-existing suites to "N", adds a "CreateContacts" / "Y"
# row to test_suite, and updates the active sheet/selection state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("test_suite")
$ws2 = $wb.Worksheets.Item("AddCustomerTest")

# --- Add the new worksheet at the end of the workbook ---------------------
$wsNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "CreateContacts"

# --- Populate cells in the exact order the new unique strings are first
#     introduced, so shared-string indices line up with the target file ---

# 1) "companyname" (new shared string)
$wsNew.Range("C1").Value = "companyname"

# 2) "CreateContacts" (new shared string) + the runmode flag for it
$ws1.Range("A5").Value = "CreateContacts"
$ws1.Range("B5").Value = "Y"

# 3) "Sumitra" (new shared string)
$wsNew.Range("A3").Value = "Sumitra"

# 4) "Dassault Systems 2" (new shared string)
$wsNew.Range("C3").Value = "Dassault Systems 2"

# 5) "Dassault Systems 1" (new shared string)
$wsNew.Range("C2").Value = "Dassault Systems 1"

# --- Remaining cells (all re-use already-existing shared strings) ---------
$ws1.Range("B2").Value = "N"
$ws1.Range("B3").Value = "N"
$ws1.Range("B4").Value = "N"

$wsNew.Range("A1").Value = "firstname"
$wsNew.Range("B1").Value = "lastname"
$wsNew.Range("D1").Value = "runmode"

$wsNew.Range("A2").Value = "Rohan"
$wsNew.Range("B2").Value = "Bhor"
$wsNew.Range("D2").Value = "Y"

$wsNew.Range("B3").Value = "Bhor"
$wsNew.Range("D3").Value = "N"

# --- Column C ("companyname") is widest on the new sheet - fit it ---------
$wsNew.Columns.Item(3).AutoFit() | Out-Null

# --- Selections / active sheet (applied last-to-first so the final
#     selected range below leaves "CreateContacts" as the active tab) ------
$ws1.Range("E21").Select() | Out-Null
$ws2.Range("E1:E2").Select() | Out-Null
$wsNew.Range("G10").Select() | Out-Null

Write-Output "CreateContacts sheet added and test_suite updated"
